$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 10183.833
$ws.Range("I2").Value = 1975.625
$ws.Range("K2").Value = 1975.625
$ws.Range("M2").Value = -1862.625
$ws.Range("H8").Value = 7147.778
$ws.Range("I8").Value = 8221.833000000001
$ws.Range("K8").Value = 24665.499
$ws.Range("M8").Value = -24526.499
$ws.Range("H18").Value = 1481.7778
$ws.Range("I18").Value = 1379.5
$ws.Range("K18").Value = 1379.5
$ws.Range("M18").Value = -1095.5
$ws.Range("H29").Value = 7858.25
$ws.Range("J29").Value = 8323.117
$ws.Range("L29").Value = 24969.351
$ws.Range("N29").Value = -25531.351
$ws.Range("H76").Value = 14291685
$ws.Range("I76").Value = 16672308
$ws.Range("J76").Value = 7950
$ws.Range("K76").Value = 16672308
$ws.Range("L76").Value = 7950
$ws.Range("M76").Value = -16671993
$ws.Range("N76").Value = -8580
$ws.Range("H79").Value = 14291685
$ws.Range("I79").Value = 16672308
$ws.Range("J79").Value = 7950
$ws.Range("K79").Value = 16672308
$ws.Range("L79").Value = 7950
$ws.Range("M79").Value = -16671216
$ws.Range("N79").Value = -10134
$ws.Range("H92").Value = 2008
$ws.Range("I92").Value = 2115.625
$ws.Range("J92").Value = 1147
$ws.Range("K92").Value = 2115.625
$ws.Range("L92").Value = 1147
$ws.Range("M92").Value = -867.625
$ws.Range("N92").Value = -3643

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2738.8
$ws.Range("I3").Value = 625
$ws.Range("J3").Value = 4148
$ws.Range("K3").Value = 625
$ws.Range("L3").Value = 4148
$ws.Range("M3").Value = -510
$ws.Range("N3").Value = -4378
$ws.Range("H45").Value = 5419.92
$ws.Range("I45").Value = 6693.5557
$ws.Range("J45").Value = 2144.8572
$ws.Range("K45").Value = 6693.5557
$ws.Range("L45").Value = 2144.8572
$ws.Range("M45").Value = -6316.5557
$ws.Range("N45").Value = -2898.8572
$ws.Range("H92").Value = 3363480
$ws.Range("J92").Value = 3363480
$ws.Range("L92").Value = 3363480
$ws.Range("N92").Value = -3368472

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25826.938
$ws.Range("I82").Value = 15548.8
$ws.Range("K82").Value = 15548.8
$ws.Range("M82").Value = -15165.8
$ws.Range("H85").Value = 25826.938
$ws.Range("I85").Value = 15548.8
$ws.Range("K85").Value = 15548.8
$ws.Range("M85").Value = -14222.8
$ws.Range("H134").Value = 7065.6665
$ws.Range("I134").Value = 3058.4285
$ws.Range("K134").Value = 9175.2855
$ws.Range("M134").Value = -6640.2855
$ws.Range("H135").Value = 97466
$ws.Range("J135").Value = 97466
$ws.Range("L135").Value = 97466
$ws.Range("N135").Value = -107606
$ws.Range("H137").Value = 127774.5
$ws.Range("J137").Value = 127774.5
$ws.Range("L137").Value = 127774.5
$ws.Range("N137").Value = -137974.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 169.875
$ws.Range("I10").Value = 169.875
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 169.875
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -30.875
$ws.Range("N10").ClearContents()
$ws.Range("H122").Value = 78231.234
$ws.Range("I122").Value = 112111.89
$ws.Range("K122").Value = 336335.67
$ws.Range("M122").Value = -333885.67

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 91050.91
$ws.Range("I23").Value = 333395
$ws.Range("J23").Value = 171.875
$ws.Range("K23").Value = 1000185
$ws.Range("L23").Value = 515.625
$ws.Range("M23").Value = -999950
$ws.Range("N23").Value = -985.625
$ws.Range("H46").Value = 1579
$ws.Range("J46").Value = 1579
$ws.Range("L46").Value = 4737
$ws.Range("N46").Value = -4919
$ws.Range("H122").Value = 18102
$ws.Range("I122").Value = 18102
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 162918
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -160468
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 35715508
$ws.Range("I139").Value = 35715508
$ws.Range("K139").Value = 107146524
$ws.Range("M139").Value = -107141384

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9342.208000000001
$ws.Range("I113").Value = 3229.1052
$ws.Range("K113").Value = 3229.1052
$ws.Range("M113").Value = -1059.1052
$ws.Range("H122").Value = 6307.4287
$ws.Range("I122").Value = 7365.1177
$ws.Range("J122").Value = 1812.25
$ws.Range("K122").Value = 22095.3531
$ws.Range("L122").Value = 5436.75
$ws.Range("M122").Value = -19645.3531
$ws.Range("N122").Value = -10336.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4341.2
$ws.Range("J22").Value = 5768
$ws.Range("L22").Value = 5768
$ws.Range("N22").Value = -6358
$ws.Range("H27").Value = 4341.2
$ws.Range("J27").Value = 5768
$ws.Range("L27").Value = 5768
$ws.Range("N27").Value = -5982
$ws.Range("H46").Value = 5821.4062
$ws.Range("J46").Value = 6109.467
$ws.Range("L46").Value = 6109.467
$ws.Range("N46").Value = -6485.467
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H101").Value = 36200
$ws.Range("J101").Value = 36200
$ws.Range("L101").Value = 36200
$ws.Range("N101").Value = -42690
$ws.Range("H133").Value = 65998.60000000001
$ws.Range("J133").Value = 65998.60000000001
$ws.Range("L133").Value = 65998.60000000001
$ws.Range("N133").Value = -71058.60000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1200
$ws.Range("J6").Value = 1200
$ws.Range("L6").Value = 1200
$ws.Range("N6").Value = -1430
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H103").Value = 68888
$ws.Range("J103").Value = 68888
$ws.Range("L103").Value = 68888
$ws.Range("N103").Value = -71232
